# "update preg visit section"
#
# Resizes/repositions the big background rectangle, the folded-corner
# AutoShape, the Composition bracket rectangle and several "entry: ..."
# labels that trail further down the page to make room for a new
# "Observation (Reference)" label that gets inserted right under the
# "Encounter (Reference)" label in the Pregnancy Progress section; also
# grows the "Pregnancy Progress (Section)" box to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

$EMU = 12700.0

# --- Rectangle 29 (id 30): big background rectangle behind the bundle ---
$shp30 = Get-ShapeById $s 30
$shp30.Left   = -107344 / $EMU
$shp30.Top    = 6571128 / $EMU
$shp30.Width  = 5764697 / $EMU
$shp30.Height = 5522806 / $EMU

# --- AutoShape 29 (id 4): folded-corner background shape, grows taller ---
$shp4 = Get-ShapeById $s 4
$shp4.Height = 11619966 / $EMU

# --- Rectangle 5 (id 5): "entry: Composition" bracket, grows taller ---
$shp5 = Get-ShapeById $s 5
$shp5.Height = 8848226 / $EMU

# --- "entry: Patient" (id 19) moves further down ---
$shp19 = Get-ShapeById $s 19
$shp19.Top = 9874018 / $EMU

# --- "entry: Encounter" (id 20) moves further down ---
$shp20 = Get-ShapeById $s 20
$shp20.Top = 10317260 / $EMU

# --- "entry: Observation" (id 15) moves further down ---
$shp15 = Get-ShapeById $s 15
$shp15.Top = 10750085 / $EMU

# --- "entry: ..." (id 18) moves further down ---
$shp18 = Get-ShapeById $s 18
$shp18.Top = 11182911 / $EMU

# --- "Pregnancy Progress (Section)" (id 28) grows taller ---
$shp28 = Get-ShapeById $s 28
$shp28.Top    = 8316926 / $EMU
$shp28.Height = 1229777 / $EMU

# --- Add new "Observation (Reference)" label under "Encounter (Reference)" ---
# Duplicate a sibling label that already has the exact same style/run
# layout (bold run + normal " (Reference)" run) and retarget it.
$shp26 = Get-ShapeById $s 26
$dup = $shp26.Duplicate()
$newShape = $dup.Item(1)

$tr = $newShape.TextFrame.TextRange
$tr.Text = "Observation (Reference)"
$tr.Characters(1, 11).Font.Bold = $true

$newShape.Left   = 757646 / $EMU
$newShape.Top    = 9069921 / $EMU
$newShape.Width  = 3845730 / $EMU
$newShape.Height = 338554 / $EMU
